# Add season-record columns (Wins / Losses / Ties) to the right of the
# existing stats table, matching the style of the other header cells,
# and fill every player row with the team's 2019 season record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns AD (30), AE (31), AF (32) are the first three empty columns
# after the existing data (which ends at AC = 29).
$colWins   = 30
$colLosses = 31
$colTies   = 32

# --- Header row (row 1) ------------------------------------------------
# Copy the existing header style (from A1, which already carries the
# bold/centered/bordered "s=1" style) onto the new header cells, then
# set their text.
$ws.Cells.Item(1, 1).Copy()
$ws.Cells.Item(1, $colWins).PasteSpecial(-4122)
$ws.Cells.Item(1, $colLosses).PasteSpecial(-4122)
$ws.Cells.Item(1, $colTies).PasteSpecial(-4122)

$ws.Cells.Item(1, $colWins).Value = "Wins"
$ws.Cells.Item(1, $colLosses).Value = "Losses"
$ws.Cells.Item(1, $colTies).Value = "Ties"

# --- Data rows (2-59) ---------------------------------------------------
# Every row corresponds to a player on the same team, so the season
# record (96 wins, 66 losses, 0 ties) is identical across all of them.
$lastRow = 59
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $colWins).Value = 96
    $ws.Cells.Item($r, $colLosses).Value = 66
    $ws.Cells.Item($r, $colTies).Value = 0
}
